$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.850.28'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.28'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.59'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3886'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3808'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.61'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.348'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.003'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08469'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.02'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.052'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.078'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001313'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.653.44'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.17'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07006'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.63'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.966'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.75'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.843.24'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.431'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.962'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.07'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.89'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.416'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.07'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.842'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.495'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.834.65'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.013'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08184'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.08%  '

$ws.Range("B36").Value = 'VeChain'

$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02902'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.11%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'

$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.608'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.80'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2674'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09154'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.58'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7572'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.423'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.46'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6936'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.454'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.098'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08288'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.15'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.222'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.93%  '
